$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Solvente" column at G (existing G:N shift right to H:O)
# and populate its header + the six data rows.
$ws.Columns("G:G").Insert()

$ws.Range("G1").Value = "Solvente"
$ws.Range("G1").Font.Bold = $true

$ws.Range("G2").Value = "Nefta Pesada Hidrotratada"
$ws.Range("G3").Value = "Gasoil Hidrotratado"
$ws.Range("G4").Value = "Agua Desmineralizada"
$ws.Range("G5").Value = "Agua Desmineralizada"
$ws.Range("G6").Value = "Reformado"

# Match the column's best-fit width from the source workbook.
$ws.Columns("G:G").ColumnWidth = 23.666666666666668
